$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "datos actualizados" timestamp banner (A1) ---
$ws.Range("A1").Value = "Datos actualizados a 13 de Mayo de 2020 a las 07:35"

# --- Hungria (row 68) ---
$ws.Range("B68").Value = 3341
$ws.Range("C68").Value = 28
$ws.Range("D68").Value = 1102
$ws.Range("E68").Value = 1809
$ws.Range("F68").Value = 45
$ws.Range("G68").Value = 5
$ws.Range("H68").Value = 430

# --- Uzbekistan (row 75) ---
$ws.Range("B75").Value = 2568
$ws.Range("C75").Value = 49
$ws.Range("D75").Value = 2040
$ws.Range("E75").Value = 518
$ws.Range("F75").Value = 8
$ws.Range("G75").Value = 0
$ws.Range("H75").Value = 10

# --- Bulgaria (row 80) ---
$ws.Range("B80").Value = 2069
$ws.Range("C80").Value = 46
$ws.Range("D80").Value = 499
$ws.Range("E80").Value = 1474
$ws.Range("F80").Value = 51
$ws.Range("G80").Value = 1
$ws.Range("H80").Value = 96

# --- Row 99 now becomes "El Salvador" (updated figures, moved ahead of Tunez) ---
$ws.Range("A99").Value = "El Salvador"
$ws.Range("B99").Value = 1037
$ws.Range("C99").Value = 39
$ws.Range("D99").Value = 374
$ws.Range("E99").Value = 643
$ws.Range("F99").Value = 15
$ws.Range("G99").Value = 0
$ws.Range("H99").Value = 20

# --- Row 100 now becomes "Tunez" (keeps its previous figures, pushed down one row) ---
$ws.Range("A100").Value = "Tunez"
$ws.Range("B100").Value = 1032
$ws.Range("C100").Value = 0
$ws.Range("D100").Value = 740
$ws.Range("E100").Value = 247
$ws.Range("F100").Value = 5
$ws.Range("G100").Value = 0
$ws.Range("H100").Value = 45
